$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Tiny correction to the existing last row's timestamp (A38)
$ws.Range("A38").Value = 44351.82712145486

# Append the new data row (row 39)
$ws.Range("A39").Value = 44352.7936396181
$ws.Range("B39").Value = 75249
$ws.Range("C39").Value = 63432
$ws.Range("D39").Value = 3379
$ws.Range("E39").Value = 2115
$ws.Range("F39").Value = 1487
$ws.Range("G39").Value = 19889
$ws.Range("H39").Value = 1446
$ws.Range("I39").Value = 891
$ws.Range("J39").Value = 202

# Match the date/time number format used by the rest of column A
$ws.Range("A39").NumberFormat = $ws.Range("A38").NumberFormat
